$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column F: header "OSMO_DEF", matching the header style used by E1..B1
$ws.Range("F1").Value = "OSMO_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New column F data rows: empty-list placeholder values
$ws.Range("F2:F5").Value = "[]"

$wb.Save()
